$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.273.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.678.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5338"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.34%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2684"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06473"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07540"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.701.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.524"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5780"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008470"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.320.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.900"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.211"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.007"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1277"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.32%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.827"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.51%  "

$ws.Range("E27").Value = "  +0.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06472"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.382"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.91%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.321"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.34%  "

$ws.Range("E31").Value = "  +0.91%  "

$ws.Range("E32").Value = "  +1.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.667"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.033"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6178"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.30%  "

$ws.Range("E36").Value = "  +1.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.699"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.253"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.112.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01621"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.04%  "

$ws.Range("E41").Value = "  +0.97%  "

$ws.Range("E42").Value = "  +0.50%  "

$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.830.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000108"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.22%  "

$ws.Range("E46").Value = "  +1.35%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.184"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.22%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05263"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.53%  "

$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.081"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.28%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4288"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.00%  "
